$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Events")

# Update the B8 story text (existing shared string, index 20): clause reordered + trailing space kept
$ws.Cells.Item(8, 2).Value = 'Мы решили просто пробежать их. Поначалу у нас всё складывалось отлично, но из-за кромешной тьмы Виктор не увидел зомби, у которого не было ног, и запнулся об него. Помочь ему мне сразу не удалось, так как мне загородили проход к нему три зомби. '

# New rows 20-27: id (A) + story (B) pairs (adds shared strings 21-28, in sheet order)
$ws.Cells.Item(20, 1).Value = 611
$ws.Cells.Item(20, 2).Value = 'Из-за того, что я тренировался с арбалетом на зомби ранее, у меня без труда получилось убить их и  помочь Виктору скинуть зомби с себя. По итогу мы смогли пробежать через оставшихся зомби и прийти на блокпост. В нас уже было хотели стрелять, так как подумали, что мы можем быть зомби. Но мы подали жесты руками и нас определили как людей. Нас пустили на военный блокпост и повели на медицинский осмотр. '

$ws.Cells.Item(21, 1).Value = 612
$ws.Cells.Item(21, 2).Value = ' Я уже не помнил слова Виктора о том, как стрелять с арбалета. Я смог убить одного, но из-за долгой перезарядки и последующего за ним промаха на меня сзади успел накинуться зомби и повалить. Сил во мне уже не оставалось, и меня загрызли. Смерть '

$ws.Cells.Item(22, 1).Value = 621
$ws.Cells.Item(22, 2).Value = 'Из-за того, что у меня ранее был опыт в бесшумных убийствах зомби, мне не составило труда сделать это еще раз. Я убил тех трёх зомби, после чего помог Виктору. По итогу мы смогли пробежать через оставшихся зомби и прийти на блокпост. В нас уже было хотели стрелять, так как подумали, что мы можем быть зомби. Но мы подали жесты руками и нас определили как людей. Нас пустили на военный блокпост и повели на медицинский осмотр. '

$ws.Cells.Item(23, 1).Value = 622
$ws.Cells.Item(23, 2).Value = 'Я решил попробовать убить их бесшумно. Но подходя к ним, я издавал слишком много шума. Они обернулись и разом накинулись на меня, после чего загрызли. '

$ws.Cells.Item(24, 1).Value = 5111
$ws.Cells.Item(24, 2).Value = 'Мы с Виктором  забрались на небольшую газель, и мы начали отстреливаться. В боезапасе было  10 стрел и 12 арбалетных болтов. На прошлой троице зомби еще в деревне я научился правильно прицеливать и перезаряжаться. Все стрелы и болты, которые мы выпускали, попадали точно в цель. Оставшихся пятерых зомби мы смогли убить с помощью ножей. Мы спустились с машины и пошли до военного блокпоста. В нас уже было хотели стрелять, так как подумали, что мы можем быть зомби. Но мы подали жесты руками и нас определили как людей. Нас пустили на военный блокпост и повели на медицинский осмотр. '

$ws.Cells.Item(25, 1).Value = 5112
$ws.Cells.Item(25, 2).Value = 'Мы с Виктором  забрались на небольшую газель и начали отстреливаться. В боезапасе было  10 стрел и 15 арбалетных болтов. Я плохо помнил инструктаж Виктора и не смог попасть всеми болтами арбалета из-за нервозности и тремора в руках. Их оставалось около 15. Но вдруг они смогли подняться на капот, а в последующем и на крышу. Нам ничего не оставалось, кроме того, что прыгать с крыши газели. Виктор смог спрыгнуть и сделать кувырок для смягчения падения, а я, прыгая с крыши, повредил себе ногу и не смог вовремя убежать. На меня накинулись зомби и загрызли. Смерть '

$ws.Cells.Item(26, 1).Value = 5331
$ws.Cells.Item(26, 2).Value = 'Мы решили прокрасться через них по-тихому. К счастью у нас уже был опыт в убийстве зомби со спины. Нам не составило труда кого-то убить, а кого-то обойти. Так что для нас это не было сложным испытанием. Мы продолжили путь до блокпоста. В нас уже было хотели стрелять, так как подумали, что мы можем быть зомби. Но мы подали жесты руками и нас определили как людей. Нас пустили на военный блокпост и повели на медицинский осмотр.'

$ws.Cells.Item(27, 1).Value = 5332
$ws.Cells.Item(27, 2).Value = 'Мы решили прокрасться через них по-тихому. Но мы никогда не оказывались в такой ситуации, поэтому, что делать мы не знали. Нас охватила легкая паника. А тем временем зомби были все ближе и больше. Еще немного подумав, мы все же решились идти. Мы начали тихонько обходить зомби, но из-за нервности мы совершили ошибку. В темноте Виктор не заметил труп и запнулся об него. От неожиданности он охнул. Зомби заметили нас и начали окружать. Их стало еще больше, и мы не смогли выбраться. Нас загрызли. '

# New column E header ("death") in row 1 -- added last so it becomes shared string 29
$ws.Cells.Item(1, 5).Value = 'death'

# Selection / scroll position to match the final view (E1 selected, scrolled so column C is leftmost)
$ws.Range("E1").Select()
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1

# Best-effort: workbook OS window position (yWindow 1905 -> 2505)
$wb.Windows.Item(1).Top = 2505
